# Applies the 2024-01-31 cryptos-list refresh: updated prices/volumes, and a
# couple of rows (20/21 and 36/37) whose ranking order swapped since the last run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.625.71'
$ws.Range('E2').Value = '  -1.72%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.304.15'
$ws.Range('E3').Value = '  -0.07%  '

# Row 4
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.49'
$ws.Range('E5').Value = '  -1.79%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.46'
$ws.Range('E6').Value = '  -4.87%  '

# Row 7
$ws.Range('E7').Value = '  -4.23%  '

# Row 8
$ws.Range('E8').Value = '  +0.11%  '

# Row 9
$ws.Range('E9').Value = '  -4.06%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.76'
$ws.Range('E10').Value = '  -4.08%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('E11').Value = '  -2.25%  '

# Row 12
$ws.Range('E12').Value = '  +0.70%  '

# Row 13
$ws.Range('E13').Value = '  -3.02%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.662.36'
$ws.Range('E14').Value = '  +0.08%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.71'
$ws.Range('E15').Value = '  +4.33%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.296.99'
$ws.Range('E16').Value = '  -0.17%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.805'
$ws.Range('E17').Value = '  +0.37%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.572.07'
$ws.Range('E18').Value = '  -1.66%  '

# Row 19
$ws.Range('E19').Value = '  -1.75%  '

# Row 20
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.07'
$ws.Range('E20').Value = '  -1.33%  '

# Row 21
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.47'
$ws.Range('E21').Value = '  -3.58%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.92'
$ws.Range('E22').Value = '  +0.09%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.90'
$ws.Range('E23').Value = '  -2.35%  '

# Row 24
$ws.Range('E24').Value = '  -2.30%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.52'
$ws.Range('E25').Value = '  -3.51%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.97'
$ws.Range('E27').Value = '  +0.74%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.29'
$ws.Range('E28').Value = '  +1.91%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.62'
$ws.Range('E29').Value = '  -4.59%  '

# Row 30
$ws.Range('E30').Value = '  -4.06%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '164.02'
$ws.Range('E31').Value = '  +0.14%  '

# Row 32
$ws.Range('E32').Value = '  +0.00%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.01'
$ws.Range('E33').Value = '  -4.16%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.42'
$ws.Range('E34').Value = '  -5.03%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.45'
$ws.Range('E35').Value = '  -1.12%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0705'
$ws.Range('E36').Value = '  -4.27%  '

# Row 37
$ws.Range('B37').Value = 'Celestia'
$ws.Range('C37').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.77'
$ws.Range('E37').Value = '  -7.99%  '

# Row 38
$ws.Range('E38').Value = '  -4.16%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.100'
$ws.Range('E40').Value = '  -5.14%  '

# Row 41
$ws.Range('E41').Value = '  -3.45%  '

# Row 42
$ws.Range('E42').Value = '  +0.77%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.966.93'
$ws.Range('E43').Value = '  -0.91%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0281'
$ws.Range('E44').Value = '  -3.13%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.44'
$ws.Range('E45').Value = '  -2.22%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.24'
$ws.Range('E46').Value = '  +1.44%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.89'
$ws.Range('E47').Value = '  -6.09%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.30'
$ws.Range('E48').Value = '  -3.72%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.528.56'
$ws.Range('E49').Value = '  -0.08%  '

# Row 50
$ws.Range('E50').Value = '  -2.59%  '

# Row 51
$ws.Range('E51').Value = '  +0.32%  '
